$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Characters(21, 2).Text = "20"
$ws.Range("C9").Characters(27, 8).Text = "5/15/2023"
$ws.Range("C9").Characters(47, 9).Text = "5/21/2023"

# --- Data table updates ---
$ws.Range("L15").Value = -37.5
$ws.Range("N15").Value = -16.666666666666
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -46.666666666666
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -20.338983050847
$ws.Range("L16").Value = 42.424242424242
$ws.Range("M16").Value = 80.769230769230
$ws.Range("N16").Value = -85.843373493975
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("I17").Value = 51
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = 13.333333333333
$ws.Range("L17").Value = 4.081632653061
$ws.Range("M17").Value = 142.857142857143
$ws.Range("N17").Value = -21.538461538461
$ws.Range("C18").Value = 13
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 160
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 4.545454545454
$ws.Range("I18").Value = 76
$ws.Range("J18").Value = 113
$ws.Range("K18").Value = -32.743362831858
$ws.Range("L18").Value = 80.952380952380
$ws.Range("M18").Value = 5.555555555555
$ws.Range("N18").Value = -77.108433734939
$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = -21.875
$ws.Range("F19").Value = 86
$ws.Range("G19").Value = 110
$ws.Range("H19").Value = -21.818181818181
$ws.Range("I19").Value = 431
$ws.Range("J19").Value = 458
$ws.Range("K19").Value = -5.895196506550
$ws.Range("L19").Value = 77.366255144032
$ws.Range("M19").Value = 3.357314148681
$ws.Range("N19").Value = -67.739520958083
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -44.444444444444
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 20
$ws.Range("M20").Value = 118.181818181818
$ws.Range("N20").Value = -92.982456140350
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = -2.325581395348
$ws.Range("F21").Value = 134
$ws.Range("G21").Value = 166
$ws.Range("H21").Value = -19.277108433734
$ws.Range("I21").Value = 634
$ws.Range("J21").Value = 703
$ws.Range("K21").Value = -9.815078236130
$ws.Range("L21").Value = 62.982005141388
$ws.Range("M21").Value = 15.904936014625
$ws.Range("N21").Value = -73.801652892562
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 29
$ws.Range("J22").Value = 42
$ws.Range("K22").Value = -30.952380952381
$ws.Range("L22").Value = 11.538461538461
$ws.Range("M22").Value = 20.833333333333
$ws.Range("C24").Value = 78
$ws.Range("D24").Value = 75
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = 276
$ws.Range("G24").Value = 346
$ws.Range("H24").Value = -20.231213872832
$ws.Range("I24").Value = 1442
$ws.Range("J24").Value = 1468
$ws.Range("K24").Value = -1.771117166212
$ws.Range("L24").Value = 98.896551724137
$ws.Range("M24").Value = 121.165644171779
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = -6.976744186046
$ws.Range("I25").Value = 138
$ws.Range("J25").Value = 136
$ws.Range("K25").Value = 1.470588235294
$ws.Range("L25").Value = 28.971962616822
$ws.Range("M25").Value = 46.808510638297
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("L26").Value = -22.222222222222
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 32
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = -8.571428571428
$ws.Range("L27").Value = 23.076923076923
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("C30").Value = 1
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 4
$ws.Range("K30").Value = -60
